$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 13 (Leve Item ID 2144)
$ws.Range("H13").Value = 9000
$ws.Range("J13").Value = 9000
$ws.Range("L13").Value = 9000
$ws.Range("N13").Value = -9338
# Row 16 (Leve Item ID 2146)
$ws.Range("H16").Value = 4300
$ws.Range("J16").Value = 4950
$ws.Range("L16").Value = 4950
$ws.Range("N16").Value = -5410
# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 4576.7646
$ws.Range("I62").Value = 2580.5
$ws.Range("J62").Value = 7428.5713
$ws.Range("K62").Value = 2580.5
$ws.Range("L62").Value = 7428.5713
$ws.Range("M62").Value = -1956.5
$ws.Range("N62").Value = -8676.5713
# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 4576.7646
$ws.Range("I65").Value = 2580.5
$ws.Range("J65").Value = 7428.5713
$ws.Range("K65").Value = 12902.5
$ws.Range("L65").Value = 37142.85649999999
$ws.Range("M65").Value = -9782.5
$ws.Range("N65").Value = -43382.85649999999
# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 5489.1665
$ws.Range("I141").Value = 4425.909
$ws.Range("J141").Value = 6388.846
$ws.Range("K141").Value = 13277.727
$ws.Range("L141").Value = 19166.538
$ws.Range("M141").Value = -8097.726999999999
$ws.Range("N141").Value = -29526.538

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 18339.207
$ws.Range("I32").Value = 20517.834
$ws.Range("K32").Value = 20517.834
$ws.Range("M32").Value = -20230.834
# Row 57 (Leve Item ID 39767)
$ws.Range("H57").Value = 10500
$ws.Range("I57").Value = 10500
$ws.Range("K57").Value = 10500
$ws.Range("M57").Value = -10016
# Row 112 (Leve Item ID 25808)
$ws.Range("H112").Value = 38174.8
$ws.Range("J112").Value = 38174.8
$ws.Range("L112").Value = 38174.8
$ws.Range("N112").Value = -41128.8
# Row 117 (Leve Item ID 26125)
$ws.Range("H117").Value = 79800
$ws.Range("J117").Value = 79800
$ws.Range("L117").Value = 79800
$ws.Range("N117").Value = -88978
# Row 119 (Leve Item ID 26287)
$ws.Range("H119").Value = 34888
$ws.Range("J119").Value = 34888
$ws.Range("L119").Value = 34888
$ws.Range("N119").Value = -44564
# Row 124 (Leve Item ID 34252)
$ws.Range("H124").Value = 18000
$ws.Range("J124").Value = 18000
$ws.Range("L124").Value = 18000
$ws.Range("N124").Value = -27820
# Row 125 (Leve Item ID 34251)
$ws.Range("H125").Value = 35000
$ws.Range("J125").Value = 35000
$ws.Range("L125").Value = 35000
$ws.Range("N125").Value = -44840
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2652.5
$ws.Range("I132").Value = 2450.7058
$ws.Range("J132").Value = 3142.5715
$ws.Range("K132").Value = 7352.117400000001
$ws.Range("L132").Value = 9427.7145
$ws.Range("M132").Value = -4822.117400000001
$ws.Range("N132").Value = -14487.7145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 11 (Leve Item ID 2481)
$ws.Range("H11").Value = 11332.667
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 11332.667
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 11332.667
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -11612.667
# Row 15 (Leve Item ID 1605)
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1209.75
$ws.Range("I99").Value = 1155.6
$ws.Range("J99").Value = 1300
$ws.Range("K99").Value = 1155.6
$ws.Range("L99").Value = 1300
$ws.Range("M99").Value = 342.4000000000001
$ws.Range("N99").Value = -4296
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2740.652
$ws.Range("I134").Value = 2768.6667
$ws.Range("J134").Value = 2639.8
$ws.Range("K134").Value = 8306.000100000001
$ws.Range("L134").Value = 7919.400000000001
$ws.Range("M134").Value = -5771.000100000001
$ws.Range("N134").Value = -12989.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 10 (Leve Item ID 1997)
$ws.Range("H10").Value = 36766.668
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 54900
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 54900
$ws.Range("M10").Value = -361
$ws.Range("N10").Value = -55178
# Row 14 (Leve Item ID 1998)
$ws.Range("H14").Value = 700
$ws.Range("I14").Value = 700
$ws.Range("K14").Value = 700
$ws.Range("M14").Value = -530
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 1570692.6
$ws.Range("I58").Value = 2021467.9
$ws.Range("J58").Value = 10316.615
$ws.Range("K58").Value = 2021467.9
$ws.Range("L58").Value = 10316.615
$ws.Range("M58").Value = -2021264.9
$ws.Range("N58").Value = -10722.615
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 1303.4286
$ws.Range("I99").Value = 1205.3334
$ws.Range("J99").Value = 1480
$ws.Range("K99").Value = 1205.3334
$ws.Range("L99").Value = 1480
$ws.Range("M99").Value = 292.6666
$ws.Range("N99").Value = -4476
# Row 125 (Leve Item ID 34297)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = ""
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 1303.4286
$ws.Range("I126").Value = 1205.3334
$ws.Range("J126").Value = 1480
$ws.Range("K126").Value = 3616.0002
$ws.Range("L126").Value = 4440
$ws.Range("M126").Value = -1146.0002
$ws.Range("N126").Value = -9380
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2622.5122
$ws.Range("I132").Value = 2242.742
$ws.Range("K132").Value = 6728.226000000001
$ws.Range("M132").Value = -4198.226000000001
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1570692.6
$ws.Range("I136").Value = 2021467.9
$ws.Range("J136").Value = 10316.615
$ws.Range("K136").Value = 6064403.699999999
$ws.Range("L136").Value = 30949.845
$ws.Range("M136").Value = -6061853.699999999
$ws.Range("N136").Value = -36049.845

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 12828720
$ws.Range("J5").Value = 23824304
$ws.Range("L5").Value = 71472912
$ws.Range("N5").Value = -71473136
# Row 15 (Leve Item ID 4661)
$ws.Range("H15").Value = 427.79166
$ws.Range("I15").Value = 165.38461
$ws.Range("J15").Value = 737.9091
$ws.Range("K15").Value = 496.15383
$ws.Range("L15").Value = 2213.7273
$ws.Range("M15").Value = -356.15383
$ws.Range("N15").Value = -2493.7273
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 888.4091
$ws.Range("I122").Value = 586.6667
$ws.Range("J122").Value = 936.0526
$ws.Range("K122").Value = 5280.0003
$ws.Range("L122").Value = 8424.473399999999
$ws.Range("M122").Value = -2830.0003
$ws.Range("N122").Value = -13324.4734
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 12828720
$ws.Range("J135").Value = 23824304
$ws.Range("L135").Value = 214418736
$ws.Range("N135").Value = -214423806

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 13 (Leve Item ID 2443)
$ws.Range("H13").Value = 1843.1818
$ws.Range("I13").Value = 1385.1428
$ws.Range("J13").Value = 2644.75
$ws.Range("K13").Value = 1385.1428
$ws.Range("L13").Value = 2644.75
$ws.Range("M13").Value = -1246.1428
$ws.Range("N13").Value = -2922.75
# Row 17 (Leve Item ID 2445)
$ws.Range("H17").Value = 9342.333000000001
$ws.Range("J17").Value = 9342.333000000001
$ws.Range("L17").Value = 9342.333000000001
$ws.Range("N17").Value = -9678.333000000001
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 3557.1428
$ws.Range("I113").Value = 3860
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 3860
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -1690
$ws.Range("N113").Value = -7140
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3112.1
$ws.Range("I132").Value = 3552.9
$ws.Range("J132").Value = 2671.3
$ws.Range("K132").Value = 10658.7
$ws.Range("L132").Value = 8013.900000000001
$ws.Range("M132").Value = -8128.700000000001
$ws.Range("N132").Value = -13073.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2910.0715
$ws.Range("I7").Value = 2333.8
$ws.Range("K7").Value = 2333.8
$ws.Range("M7").Value = -2221.8
# Row 110 (Leve Item ID 25809)
$ws.Range("H110").Value = 68548
$ws.Range("J110").Value = 68548
$ws.Range("L110").Value = 68548
$ws.Range("N110").Value = -76728
# Row 114 (Leve Item ID 25990)
$ws.Range("H114").Value = 70200
$ws.Range("J114").Value = 70200
$ws.Range("L114").Value = 70200
$ws.Range("N114").Value = -78878
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2910.0715
$ws.Range("I126").Value = 2333.8
$ws.Range("K126").Value = 7001.400000000001
$ws.Range("M126").Value = -4531.400000000001
# Row 127 (Leve Item ID 34401)
$ws.Range("H127").Value = 60529
$ws.Range("J127").Value = 60529
$ws.Range("L127").Value = 60529
$ws.Range("N127").Value = -70449
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2662.7632
$ws.Range("I132").Value = 2168.1875
$ws.Range("J132").Value = 5300.5
$ws.Range("K132").Value = 6504.5625
$ws.Range("L132").Value = 15901.5
$ws.Range("M132").Value = -3974.5625
$ws.Range("N132").Value = -20961.5
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 3963.6965
$ws.Range("I136").Value = 2440.4333
$ws.Range("J136").Value = 5721.3076
$ws.Range("K136").Value = 7321.2999
$ws.Range("L136").Value = 17163.9228
$ws.Range("M136").Value = -4771.2999
$ws.Range("N136").Value = -22263.9228

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 3 (Leve Item ID 3309)
$ws.Range("H3").Value = 2502500
$ws.Range("I3").Value = 2502500
$ws.Range("K3").Value = 2502500
$ws.Range("M3").Value = -2502386
# Row 5 (Leve Item ID 3515)
$ws.Range("H5").Value = 13339963
$ws.Range("J5").Value = 13339963
$ws.Range("L5").Value = 13339963
$ws.Range("N5").Value = -13340187
# Row 8 (Leve Item ID 2999)
$ws.Range("H8").Value = 4000
$ws.Range("J8").Value = 4000
$ws.Range("L8").Value = 4000
$ws.Range("N8").Value = -4280
# Row 119 (Leve Item ID 26289)
$ws.Range("H119").Value = 79800
$ws.Range("J119").Value = 79800
$ws.Range("L119").Value = 79800
$ws.Range("N119").Value = -89476
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 1275.9333
$ws.Range("I126").Value = 1232.5186
$ws.Range("J126").Value = 1666.6666
$ws.Range("K126").Value = 3697.5558
$ws.Range("L126").Value = 4999.9998
$ws.Range("M126").Value = -1227.5558
$ws.Range("N126").Value = -9939.9998
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2535.7273
$ws.Range("I132").Value = 1448.3889
$ws.Range("J132").Value = 3840.5334
$ws.Range("K132").Value = 4345.1667
$ws.Range("L132").Value = 11521.6002
$ws.Range("M132").Value = -1815.1667
$ws.Range("N132").Value = -16581.6002
